$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F header "time_taken", copying E1's format (bold/border/alignment)
# so the new header cell reuses the same style as the other header cells.
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Populate time_taken values for each data row (plain, unstyled like the other data cells)
$ws.Range("F2").Value = "2021-10-05 13:38:49.882856"
$ws.Range("F3").Value = "2021-10-05 13:38:49.882864"
$ws.Range("F4").Value = "2021-10-05 13:38:49.882866"
$ws.Range("F5").Value = "2021-10-05 13:38:49.882868"
$ws.Range("F6").Value = "2021-10-05 13:38:49.882871"
$ws.Range("F7").Value = "2021-10-05 13:38:49.882873"
$ws.Range("F8").Value = "2021-10-05 13:38:49.882874"
